$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1684, 1).Value = "2025-06-09 00:19:24"
$ws.Cells.Item(1684, 2).Value = 515
$ws.Cells.Item(1684, 3).Value = 21.02
$ws.Cells.Item(1684, 4).Value = 115.92
$ws.Cells.Item(1684, 5).Value = 26.6
$ws.Cells.Item(1684, 6).Value = 67.7

$ws.Cells.Item(1685, 1).Value = "2025-06-09 00:19:25"
$ws.Cells.Item(1685, 2).Value = 516
$ws.Cells.Item(1685, 3).Value = 21.02
$ws.Cells.Item(1685, 4).Value = 115.87
$ws.Cells.Item(1685, 5).Value = 26.6
$ws.Cells.Item(1685, 6).Value = 67.7

$ws.Cells.Item(1686, 1).Value = "2025-06-09 00:19:26"
$ws.Cells.Item(1686, 2).Value = 516
$ws.Cells.Item(1686, 3).Value = 21.02
$ws.Cells.Item(1686, 4).Value = 115.75
$ws.Cells.Item(1686, 5).Value = 26.6
$ws.Cells.Item(1686, 6).Value = 67.7

$ws.Cells.Item(1687, 1).Value = "2025-06-09 00:19:27"
$ws.Cells.Item(1687, 2).Value = 516
$ws.Cells.Item(1687, 3).Value = 20.53
$ws.Cells.Item(1687, 4).Value = 115.75
$ws.Cells.Item(1687, 5).Value = 26.6
$ws.Cells.Item(1687, 6).Value = 67.7

$ws.Cells.Item(1688, 1).Value = "2025-06-09 00:19:29"
$ws.Cells.Item(1688, 2).Value = 516
$ws.Cells.Item(1688, 3).Value = 20.53
$ws.Cells.Item(1688, 4).Value = 115.84
$ws.Cells.Item(1688, 5).Value = 26.6
$ws.Cells.Item(1688, 6).Value = 67.7

$ws.Cells.Item(1689, 1).Value = "2025-06-09 00:19:30"
$ws.Cells.Item(1689, 2).Value = 516
$ws.Cells.Item(1689, 3).Value = 21.51
$ws.Cells.Item(1689, 4).Value = 115.77
$ws.Cells.Item(1689, 5).Value = 26.6
$ws.Cells.Item(1689, 6).Value = 67.7

$ws.Cells.Item(1690, 1).Value = "2025-06-09 00:19:31"
$ws.Cells.Item(1690, 2).Value = 516
$ws.Cells.Item(1690, 3).Value = 21.51
$ws.Cells.Item(1690, 4).Value = 69.63
$ws.Cells.Item(1690, 5).Value = 26.6
$ws.Cells.Item(1690, 6).Value = 67.7

$ws.Cells.Item(1691, 1).Value = "2025-06-09 00:19:37"
$ws.Cells.Item(1691, 2).Value = 516
$ws.Cells.Item(1691, 3).Value = 21.02
$ws.Cells.Item(1691, 4).Value = 70.52
$ws.Cells.Item(1691, 5).Value = 26.6
$ws.Cells.Item(1691, 6).Value = 67.7

$ws.Cells.Item(1692, 1).Value = "2025-06-09 00:19:43"
$ws.Cells.Item(1692, 2).Value = 516
$ws.Cells.Item(1692, 3).Value = 21.51
$ws.Cells.Item(1692, 4).Value = 70.52
$ws.Cells.Item(1692, 5).Value = 26.6
$ws.Cells.Item(1692, 6).Value = 67.7

$ws.Cells.Item(1693, 1).Value = "2025-06-09 00:19:44"
$ws.Cells.Item(1693, 2).Value = 516
$ws.Cells.Item(1693, 3).Value = 21.02
$ws.Cells.Item(1693, 4).Value = 69.7
$ws.Cells.Item(1693, 5).Value = 26.6
$ws.Cells.Item(1693, 6).Value = 67.7

$ws.Cells.Item(1694, 1).Value = "2025-06-09 00:19:45"
$ws.Cells.Item(1694, 2).Value = 516
$ws.Cells.Item(1694, 3).Value = 21.51
$ws.Cells.Item(1694, 4).Value = 69.27
$ws.Cells.Item(1694, 5).Value = 26.6
$ws.Cells.Item(1694, 6).Value = 67.7

$ws.Cells.Item(1695, 1).Value = "2025-06-09 00:19:46"
$ws.Cells.Item(1695, 2).Value = 515
$ws.Cells.Item(1695, 3).Value = 21.02
$ws.Cells.Item(1695, 4).Value = 69.7
$ws.Cells.Item(1695, 5).Value = 26.6
$ws.Cells.Item(1695, 6).Value = 67.7

$ws.Cells.Item(1696, 1).Value = "2025-06-09 00:19:47"
$ws.Cells.Item(1696, 2).Value = 515
$ws.Cells.Item(1696, 3).Value = 21.51
$ws.Cells.Item(1696, 4).Value = 69.72
$ws.Cells.Item(1696, 5).Value = 26.6
$ws.Cells.Item(1696, 6).Value = 67.7

$ws.Cells.Item(1697, 1).Value = "2025-06-09 00:19:48"
$ws.Cells.Item(1697, 2).Value = 515
$ws.Cells.Item(1697, 3).Value = 21.51
$ws.Cells.Item(1697, 4).Value = 70.52
$ws.Cells.Item(1697, 5).Value = 26.6
$ws.Cells.Item(1697, 6).Value = 67.7

$ws.Cells.Item(1698, 1).Value = "2025-06-09 00:19:49"
$ws.Cells.Item(1698, 2).Value = 515
$ws.Cells.Item(1698, 3).Value = 21.51
$ws.Cells.Item(1698, 4).Value = 70.06999999999999
$ws.Cells.Item(1698, 5).Value = 26.6
$ws.Cells.Item(1698, 6).Value = 67.7

$ws.Cells.Item(1699, 1).Value = "2025-06-09 00:19:50"
$ws.Cells.Item(1699, 2).Value = 515
$ws.Cells.Item(1699, 3).Value = 21.02
$ws.Cells.Item(1699, 4).Value = 70.06
$ws.Cells.Item(1699, 5).Value = 26.6
$ws.Cells.Item(1699, 6).Value = 67.7

$ws.Cells.Item(1700, 1).Value = "2025-06-09 00:19:52"
$ws.Cells.Item(1700, 2).Value = 515
$ws.Cells.Item(1700, 3).Value = 21.51
$ws.Cells.Item(1700, 4).Value = 69.7
$ws.Cells.Item(1700, 5).Value = 26.7
$ws.Cells.Item(1700, 6).Value = 67.7

$ws.Cells.Item(1701, 1).Value = "2025-06-09 00:19:53"
$ws.Cells.Item(1701, 2).Value = 515
$ws.Cells.Item(1701, 3).Value = 21.51
$ws.Cells.Item(1701, 4).Value = 69.26000000000001
$ws.Cells.Item(1701, 5).Value = 26.7
$ws.Cells.Item(1701, 6).Value = 67.7

$ws.Cells.Item(1702, 1).Value = "2025-06-09 00:19:54"
$ws.Cells.Item(1702, 2).Value = 516
$ws.Cells.Item(1702, 3).Value = 45.45
$ws.Cells.Item(1702, 4).Value = 69.63
$ws.Cells.Item(1702, 5).Value = 26.6
$ws.Cells.Item(1702, 6).Value = 67.59999999999999

$ws.Cells.Item(1703, 1).Value = "2025-06-09 00:19:55"
$ws.Cells.Item(1703, 2).Value = 516
$ws.Cells.Item(1703, 3).Value = 21.51
$ws.Cells.Item(1703, 4).Value = 69.63
$ws.Cells.Item(1703, 5).Value = 26.6
$ws.Cells.Item(1703, 6).Value = 67.59999999999999

$ws.Cells.Item(1704, 1).Value = "2025-06-09 00:19:56"
$ws.Cells.Item(1704, 2).Value = 516
$ws.Cells.Item(1704, 3).Value = 21.51
$ws.Cells.Item(1704, 4).Value = 68.83
$ws.Cells.Item(1704, 5).Value = 26.6
$ws.Cells.Item(1704, 6).Value = 67.59999999999999

$ws.Cells.Item(1705, 1).Value = "2025-06-09 00:19:57"
$ws.Cells.Item(1705, 2).Value = 516
$ws.Cells.Item(1705, 3).Value = 21.51
$ws.Cells.Item(1705, 4).Value = 69.27
$ws.Cells.Item(1705, 5).Value = 26.6
$ws.Cells.Item(1705, 6).Value = 67.59999999999999

$ws.Cells.Item(1706, 1).Value = "2025-06-09 00:19:58"
$ws.Cells.Item(1706, 2).Value = 516
$ws.Cells.Item(1706, 3).Value = 20.53
$ws.Cells.Item(1706, 4).Value = 69.27
$ws.Cells.Item(1706, 5).Value = 26.7
$ws.Cells.Item(1706, 6).Value = 67.59999999999999

$ws.Cells.Item(1707, 1).Value = "2025-06-09 00:19:59"
$ws.Cells.Item(1707, 2).Value = 516
$ws.Cells.Item(1707, 3).Value = 21.51
$ws.Cells.Item(1707, 4).Value = 68.83
$ws.Cells.Item(1707, 5).Value = 26.7
$ws.Cells.Item(1707, 6).Value = 67.59999999999999

$ws.Cells.Item(1708, 1).Value = "2025-06-09 00:20:00"
$ws.Cells.Item(1708, 2).Value = 516
$ws.Cells.Item(1708, 3).Value = 20.53
$ws.Cells.Item(1708, 4).Value = 69.27
$ws.Cells.Item(1708, 5).Value = 26.7
$ws.Cells.Item(1708, 6).Value = 67.59999999999999

$ws.Cells.Item(1709, 1).Value = "2025-06-09 00:20:01"
$ws.Cells.Item(1709, 2).Value = 516
$ws.Cells.Item(1709, 3).Value = 21.02
$ws.Cells.Item(1709, 4).Value = 68.83
$ws.Cells.Item(1709, 5).Value = 26.7
$ws.Cells.Item(1709, 6).Value = 67.59999999999999
